# Add api response formatting
# Append a new row (row 93) of API response data to each of the four
# worksheets (MID_LFT_#1, MID_LFT_#2, MID_PLT_#1, MID_PLT_#2), mirroring
# the previous row's payload but with a fresh timestamp.

$wb = $excel.ActiveWorkbook

$timestamp = 45879.4658912037
$dateFormat = "YYYY-MM-DD HH:MM:SS"

function Add-ApiRow {
    param(
        [string]$SheetName,
        [string]$Total,
        [string]$Id,
        [string]$ActualLen,
        [string]$Checksum,
        [double]$TotalDec,
        [double]$IdDec,
        [double]$ActualLenDec,
        [double]$ChecksumDec
    )

    $ws = $wb.Worksheets.Item($SheetName)
    $row = 93

    $ws.Cells.Item($row, 1).Value = $timestamp
    $ws.Cells.Item($row, 1).NumberFormat = $dateFormat

    $ws.Cells.Item($row, 2).Value = $Total
    $ws.Cells.Item($row, 3).Value = $Id
    $ws.Cells.Item($row, 4).Value = $ActualLen
    $ws.Cells.Item($row, 5).Value = $Checksum
    $ws.Cells.Item($row, 6).Value = $TotalDec
    $ws.Cells.Item($row, 7).Value = $IdDec
    $ws.Cells.Item($row, 8).Value = $ActualLenDec
    $ws.Cells.Item($row, 9).Value = $ChecksumDec
}

# ID_DEC values written as plain (non-exponential) decimal literals so the
# PowerShell parser (which has no scientific-notation support) still yields
# the exact double 5.68631262647113e+23 / 5.68432987514711e+23.

Add-ApiRow `
    "MID_LFT_#1" `
    "0x01,0x90" `
    "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c," `
    "0x01,0x1C" `
    "0x07" `
    400 `
    568631262647113000000000.0 `
    284 `
    7

Add-ApiRow `
    "MID_LFT_#2" `
    "0x01,0x7c" `
    "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," `
    "0x01,0x28" `
    "0x19" `
    380 `
    568432987514711000000000.0 `
    296 `
    25

Add-ApiRow `
    "MID_PLT_#1" `
    "0x00,0x6e" `
    "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," `
    "0x00,0x5E" `
    "0x15" `
    110 `
    568631262647113000000000.0 `
    94 `
    15

Add-ApiRow `
    "MID_PLT_#2" `
    "0x00,0x82" `
    "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," `
    "0x00,0x75" `
    "0x9" `
    130 `
    568631262647113000000000.0 `
    117 `
    9
